$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in row 4 (trialTrain = 3)
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("H4").Value = 46

# Move the selection to C4, like the saved workbook state
$ws.Range("C4").Select()
